# Implements "proper handling of invalid submits (empty outputs)":
# appends 5 new log rows (rows 4-8) to the "Data" worksheet, mirroring
# the existing log-row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting the
# engine's smart "looks like a date" auto-detection (e.g. for strings
# such as "2018.03.05") turn it into a serial date number/format.
# We do this by stashing the literal text as the cached string-result of
# a temporary formula in a scratch cell far away from the used range,
# then copying that *value* (PasteSpecial values-only) into the target
# cell. A formula's cached text result pastes back as plain text (shared
# string), bypassing the "new user input" date/number parser entirely
# and leaving no trace in the styles part.
$scratch = $ws.Range("ZZ1000")
function Set-TextValue {
    param($cellAddr, [string]$text)
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

function Set-Row {
    param($r, $dateText, $timeText, $neuron, $d, $e, $f, $g, $h, $i, $j, $k, $l)

    Set-TextValue "A$r" $dateText
    Set-TextValue "B$r" $timeText
    Set-TextValue "C$r" $neuron

    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e

    if ($f -is [string]) { Set-TextValue "F$r" $f } else { $ws.Range("F$r").Value = $f }

    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $i
    $ws.Range("J$r").Value = $j

    if ($k -is [string]) { Set-TextValue "K$r" $k } else { $ws.Range("K$r").Value = $k }
    if ($l -is [string]) { Set-TextValue "L$r" $l } else { $ws.Range("L$r").Value = $l }
}

Set-Row 4 "2018.03.05" "11:57:07" "RS" 10   50 0.1    0.96 3495 0.35 1 32.62626262626263    "N/A"
Set-Row 5 "2018.03.05" "12:00:40" "RS" 1000 50 "N/A"  1    3499 0.18 1 0.7014028056112225   "N/A"
Set-Row 6 "2018.03.05" "13:44:08" "RS" 10   50 "N/A"  1    3499 0.2  1 31.1623246492986     "N/A"
Set-Row 7 "2018.03.05" "13:48:20" "RS" 10   50 "N/A"  1    3499 0.2  1 31.1623246492986     "N/A"
Set-Row 8 "2018.03.05" "14:07:39" "RS" 10   1  "N/A"  1    3499 0.06 0 "N/A"                "N/A"

# Clean up the scratch cell so it leaves no residue in the saved sheet.
$scratch.Clear()
